$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Fisher_2009_tableS3")
$ws.Activate()

$ws.Range("E8").Value = "n"

$ws.Range("D9").Value = "exposed"
$ws.Range("E9").Formula = "=SUM(E2:E3)"

$ws.Range("D10").Value = "semi-exposed"
$ws.Range("E10").Formula = "=SUM(E4:E5)"

$ws.Range("D11").Value = "sheltered"
$ws.Range("E11").Formula = "=SUM(E6:E7)"

$null = $ws.Range("E9").Select()
